# Implementa la "plantilla semanal (max 4 cajeros)":
# Reorganiza los cajeros/horarios del turno de tarde-noche en la hoja de
# ubicacion diaria de cajeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Caja 3 (filas 8-9)
$ws.Range("F8").Value = "ERIQUE CALLE, MARIA ANTONIETA"
$ws.Range("G8").Value = "15:45-19:15"
$ws.Range("B9").Value = "YANQUI BRAVO, MIRIAN LUZ"
$ws.Range("C9").Value = "12:00-15:45"
$ws.Range("F9").Value = "VILCAPOMA CHILIN, JULISSA JAZMIN"
$ws.Range("G9").Value = "19:15-22:00"

# Caja 5 (fila 12) - se vacia
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""

# Caja 6 (fila 13)
$ws.Range("B13").Value = "HEREDIA CAHUAYA, SUSAN NAYELLI"
$ws.Range("C13").Value = "13:45-22:45"

# Fila 15
$ws.Range("B15").Value = "HURTADO SAMPLINI, JACK FERNANDO"
$ws.Range("C15").Value = "14:00-21:30"

# Caja 7 (fila 16)
$ws.Range("F16").Value = "RAMOS TINOCO, JORDAN JAVIER"
$ws.Range("G16").Value = "14:15-22:45"

# Caja 8 (fila 18)
$ws.Range("F18").Value = "LA ROSA EUSEBIO, SHADIA SHAMIRA"
$ws.Range("G18").Value = "18:15-22:00"

# Caja 9 (fila 20) - se vacia
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""

# Fila 21
$ws.Range("B21").Value = "PEREZ GORMAS, ANTHONY"
$ws.Range("C21").Value = "16:00-19:45"

# Fila 23
$ws.Range("B23").Value = "GARRIDO SOTO, VICTORIA CELESTE"
$ws.Range("C23").Value = "16:45-20:30"

# Fila 25
$ws.Range("B25").Value = "LEON TICONA, MARIA FERNANDA"
$ws.Range("C25").Value = "17:00-20:45"

# Caja 12 (fila 26)
$ws.Range("B26").Value = "DEL AGUILA MURAYARI, DARLA"
$ws.Range("C26").Value = "12:00-21:00"

# Caja 13 (fila 28)
$ws.Range("B28").Value = "IDELFONSO MOTTA, JHOSSEP ANGELO"
$ws.Range("C28").Value = "17:30-21:15"

# Caja 14 (fila 30) - hora sin cambios
$ws.Range("B30").Value = "INGA DELGADO, CARLOS DANIEL"

# Caja 15 (fila 32)
$ws.Range("B32").Value = "BRENIS LARTIGA, SEBASTIAN"
$ws.Range("C32").Value = "18:00-21:45"

# Caja 21 (fila 43)
$ws.Range("B43").Value = "MARTINEZ PAZ, ROCIO ESPERANZA"
$ws.Range("C43").Value = "09:00-18:00"

# Caja 22 (fila 45-46)
$ws.Range("B45").Value = "JIMENEZ TORERO, ASTRID GERALDINE"
$ws.Range("C45").Value = "13:15-13:30"
$ws.Range("B46").Value = "AYQUIPA MONTENEGRO, VALERIA ESTEFANY"
$ws.Range("C46").Value = "13:45-22:00"
